$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge the "TUE May 28" / " 14:29:52 IST 2019" runs into one run.
#    (Find/Replace spanning the two runs naturally collapses them into a
#    single run with the shared formatting, matching the target XML.)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("TUE May 28 14:29:52 IST 2019", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "TUE May 28 14:29:52 IST 2019", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Locate the last "CASH AND CLEARD" occurrence (the end of the last
#    purchase record in the document) and append a brand new purchase
#    record (THU May 30 / JAYAMMA / bill 13077 / CARROT + BEET) right
#    after it, before the trailing blank paragraphs.
# ---------------------------------------------------------------------
$rng = $d.Content
$lastFound = $null
while ($rng.Find.Execute("CASH AND CLEARD", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $lastFound = $d.Range($rng.Start, $rng.End)
    $rng.Collapse(0)
    $rng.End = $d.Content.End
}

$lastPara = $lastFound.Paragraphs(1)
$nextPara = $lastPara.Next()
$insertAt = $nextPara.Index

function Add-Line([int]$atIndex, [string]$text, [bool]$red) {
    $target = $d.Paragraphs($atIndex)
    $target.Range.InsertParagraphBefore()
    $newPara = $d.Paragraphs($atIndex)
    $nr = $newPara.Range
    if ($red) {
        $nr.Font.Color = 255
    }
    if ($text.Length -gt 0) {
        $nr.InsertAfter($text)
    }
    return $atIndex + 1
}

$insertAt = Add-Line $insertAt "THU May 30 13:31:01 IST 2019" $false
$insertAt = Add-Line $insertAt "Person Name`t`t`t`t- JAYAMMA" $false
$insertAt = Add-Line $insertAt "Bill number`t`t`t`t- 13077" $false
$insertAt = Add-Line $insertAt "---------------------------------------------------------------" $false
$insertAt = Add-Line $insertAt "Item Name`t`t`t`t- CARROT" $false
$insertAt = Add-Line $insertAt "Number of Pockets`t`t`t- 1" $false
$insertAt = Add-Line $insertAt "Number of KGs`t`t`t- 56" $false
$insertAt = Add-Line $insertAt "Rate`t`t`t`t`t- 40" $false
$insertAt = Add-Line $insertAt "Total Price`t`t`t`t- 2240.0" $false
$insertAt = Add-Line $insertAt "Amount balance`t`t`t- 2240.0" $false
$insertAt = Add-Line $insertAt "" $false
$insertAt = Add-Line $insertAt "Item Name`t`t`t`t- BEET" $false
$insertAt = Add-Line $insertAt "Amount Received`t`t`t- 2240" $true
$insertAt = Add-Line $insertAt "Amount Received mode`t`t- CASH AND CLEARD" $false
$insertAt = Add-Line $insertAt "" $false

Write-Host "Edit complete"
